# Updates the cryptocurrency price/volume table (cols B-E, rows 2-51) to match
# the "Updated cryptos list" GitHub Actions commit: new Price (D) / Volume(1h) (E)
# figures, a couple of corrected "-0.15%"-style signs, and the PEPE / TheGraph /
# ThetaToken rows (37-39) being re-ordered with refreshed data.
#
# Set-TextValue writes a value that Excel would otherwise auto-parse as a number
# (e.g. "3.63") by prefixing it with an apostrophe (forces text, like the source
# file's inlineStr cells) and then resets Style back to "Normal" so the cell does
# not end up pinned to a new "quote prefixed" style.
function Set-TextValue($Cell, $Text) {
    $Cell.Value = "'" + $Text
    $Cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "72.576.62"
$ws.Cells.Item(2,5).Value = "  +4.62%  "
$ws.Cells.Item(3,4).Value = "4.076.56"
$ws.Cells.Item(3,5).Value = "  +4.24%  "
$ws.Cells.Item(4,5).Value = "  -0.15%  "
Set-TextValue $ws.Cells.Item(5,4) "520.51"
$ws.Cells.Item(5,5).Value = "  -1.80%  "
Set-TextValue $ws.Cells.Item(6,4) "148.21"
$ws.Cells.Item(6,5).Value = "  +2.84%  "
Set-TextValue $ws.Cells.Item(7,4) "0.741"
$ws.Cells.Item(7,5).Value = "  +20.77%  "
Set-TextValue $ws.Cells.Item(8,4) "0.998"
$ws.Cells.Item(8,5).Value = "  -0.06%  "
Set-TextValue $ws.Cells.Item(9,4) "0.776"
$ws.Cells.Item(9,5).Value = "  +8.36%  "
Set-TextValue $ws.Cells.Item(10,4) "0.176"
$ws.Cells.Item(10,5).Value = "  +1.57%  "
$ws.Cells.Item(11,5).Value = "  -2.06%  "
Set-TextValue $ws.Cells.Item(12,4) "47.37"
$ws.Cells.Item(12,5).Value = "  +12.62%  "
Set-TextValue $ws.Cells.Item(13,4) "11.17"
$ws.Cells.Item(13,5).Value = "  +8.83%  "
$ws.Cells.Item(14,4).Value = "4.695.67"
$ws.Cells.Item(14,5).Value = "  +3.41%  "
$ws.Cells.Item(15,4).Value = "4.063.29"
$ws.Cells.Item(15,5).Value = "  +4.50%  "
Set-TextValue $ws.Cells.Item(16,4) "21.30"
$ws.Cells.Item(16,5).Value = "  +8.08%  "
Set-TextValue $ws.Cells.Item(17,4) "14.24"
$ws.Cells.Item(17,5).Value = "  +1.86%  "
Set-TextValue $ws.Cells.Item(18,4) "1.22"
$ws.Cells.Item(18,5).Value = "  -0.58%  "
$ws.Cells.Item(19,5).Value = "  -1.47%  "
$ws.Cells.Item(20,4).Value = "72.424.08"
$ws.Cells.Item(20,5).Value = "  +4.44%  "
Set-TextValue $ws.Cells.Item(21,4) "447.77"
$ws.Cells.Item(21,5).Value = "  +4.74%  "
Set-TextValue $ws.Cells.Item(22,4) "105.00"
$ws.Cells.Item(22,5).Value = "  +18.51%  "
Set-TextValue $ws.Cells.Item(23,4) "3.63"
$ws.Cells.Item(23,5).Value = "  +7.44%  "
Set-TextValue $ws.Cells.Item(24,4) "14.82"
$ws.Cells.Item(24,5).Value = "  +5.03%  "
Set-TextValue $ws.Cells.Item(25,4) "4.01"
$ws.Cells.Item(25,5).Value = "  -1.43%  "
Set-TextValue $ws.Cells.Item(26,4) "11.47"
$ws.Cells.Item(26,5).Value = "  +0.07%  "
Set-TextValue $ws.Cells.Item(27,4) "11.13"
$ws.Cells.Item(27,5).Value = "  +5.05%  "
Set-TextValue $ws.Cells.Item(28,4) "38.14"
$ws.Cells.Item(28,5).Value = "  +4.61%  "
Set-TextValue $ws.Cells.Item(29,4) "5.81"
$ws.Cells.Item(29,5).Value = "  +2.18%  "
Set-TextValue $ws.Cells.Item(30,4) "3.18"
$ws.Cells.Item(30,5).Value = "  +13.10%  "
Set-TextValue $ws.Cells.Item(31,4) "13.77"
$ws.Cells.Item(31,5).Value = "  +4.75%  "
Set-TextValue $ws.Cells.Item(32,4) "0.131"
$ws.Cells.Item(32,5).Value = "  +3.80%  "
Set-TextValue $ws.Cells.Item(33,4) "681.77"
$ws.Cells.Item(33,5).Value = "  +1.23%  "
Set-TextValue $ws.Cells.Item(34,4) "6.87"
$ws.Cells.Item(34,5).Value = "  +15.85%  "
Set-TextValue $ws.Cells.Item(35,4) "67.41"
$ws.Cells.Item(35,5).Value = "  -2.23%  "
Set-TextValue $ws.Cells.Item(36,4) "43.53"
$ws.Cells.Item(36,5).Value = "  +8.72%  "
$ws.Cells.Item(37,2).Value = "PEPE"
$ws.Cells.Item(37,3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(37,4).Value = "0.0₃0865"
$ws.Cells.Item(37,5).Value = "  -1.67%  "
$ws.Cells.Item(38,2).Value = "TheGraph"
$ws.Cells.Item(38,3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Cells.Item(38,4) "0.432"
$ws.Cells.Item(38,5).Value = "  -1.52%  "
$ws.Cells.Item(39,2).Value = "ThetaToken"
$ws.Cells.Item(39,3).Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue $ws.Cells.Item(39,4) "3.62"
$ws.Cells.Item(39,5).Value = "  +12.27%  "
Set-TextValue $ws.Cells.Item(40,4) "0.152"
$ws.Cells.Item(40,5).Value = "  +2.22%  "
$ws.Cells.Item(41,5).Value = "  +0.02%  "
$ws.Cells.Item(42,5).Value = "  +3.87%  "
Set-TextValue $ws.Cells.Item(43,4) "0.998"
$ws.Cells.Item(43,5).Value = "  -0.24%  "
Set-TextValue $ws.Cells.Item(44,4) "0.161"
$ws.Cells.Item(44,5).Value = "  +14.30%  "
Set-TextValue $ws.Cells.Item(45,4) "3.23"
$ws.Cells.Item(45,5).Value = "  +1.89%  "
$ws.Cells.Item(46,5).Value = "  -2.51%  "
Set-TextValue $ws.Cells.Item(47,4) "3.47"
$ws.Cells.Item(47,5).Value = "  +2.55%  "
$ws.Cells.Item(48,5).Value = "  +2.79%  "
Set-TextValue $ws.Cells.Item(49,4) "9.13"
$ws.Cells.Item(49,5).Value = "  +8.05%  "
Set-TextValue $ws.Cells.Item(50,4) "3.35"
$ws.Cells.Item(50,5).Value = "  +2.90%  "
Set-TextValue $ws.Cells.Item(51,4) "2.09"
$ws.Cells.Item(51,5).Value = "  +1.67%  "
